$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1861.8772
$ws.Range("I132").Value = 1349.5306
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 4048.5918
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -1518.5918
$ws.Range("N132").Value = -20060

$ws.Range("H138").Value = 2746.2727
$ws.Range("I138").Value = 1808.5405
$ws.Range("J138").Value = 3942.6897
$ws.Range("K138").Value = 5425.6215
$ws.Range("L138").Value = 11828.0691
$ws.Range("M138").Value = -285.6215000000002
$ws.Range("N138").Value = -22108.0691

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 590589.25
$ws.Range("I61").Value = 2397.5
$ws.Range("J61").Value = 1430863.1
$ws.Range("K61").Value = 2397.5
$ws.Range("L61").Value = 1430863.1
$ws.Range("M61").Value = -2185.5
$ws.Range("N61").Value = -1431287.1

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws.Range("H74").Value = 4481.6
$ws.Range("I74").Value = 1171.6296
$ws.Range("J74").Value = 34271.332
$ws.Range("K74").Value = 1171.6296
$ws.Range("L74").Value = 34271.332
$ws.Range("M74").Value = -297.6296
$ws.Range("N74").Value = -36019.332

$ws.Range("H77").Value = 4481.6
$ws.Range("I77").Value = 1171.6296
$ws.Range("J77").Value = 34271.332
$ws.Range("K77").Value = 5858.148
$ws.Range("L77").Value = 171356.66
$ws.Range("M77").Value = -1490.148
$ws.Range("N77").Value = -180092.66

$ws.Range("H136").Value = 590589.25
$ws.Range("I136").Value = 2397.5
$ws.Range("J136").Value = 1430863.1
$ws.Range("K136").Value = 7192.5
$ws.Range("L136").Value = 4292589.300000001
$ws.Range("M136").Value = -4642.5
$ws.Range("N136").Value = -4297689.300000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1312.6471
$ws.Range("I99").Value = 886.1429000000001
$ws.Range("J99").Value = 1611.2
$ws.Range("K99").Value = 886.1429000000001
$ws.Range("L99").Value = 1611.2
$ws.Range("M99").Value = 611.8570999999999
$ws.Range("N99").Value = -4607.2

$ws.Range("H105").Value = 2690
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 2690
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 2690
$ws.Range("N105").Value = -6184
$ws.Range("M105").ClearContents()

$ws.Range("H134").Value = 1736.4482
$ws.Range("I134").Value = 1190.5625
$ws.Range("J134").Value = 2408.3076
$ws.Range("K134").Value = 3571.6875
$ws.Range("L134").Value = 7224.9228
$ws.Range("M134").Value = -1036.6875
$ws.Range("N134").Value = -12294.9228

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 35715504
$ws.Range("I58").Value = 76923790
$ws.Range("J58").Value = 1652.2667
$ws.Range("K58").Value = 76923790
$ws.Range("L58").Value = 1652.2667
$ws.Range("M58").Value = -76923587
$ws.Range("N58").Value = -2058.2667

$ws.Range("H132").Value = 5347.45
$ws.Range("I132").Value = 5368.64
$ws.Range("J132").Value = 5312.1333
$ws.Range("K132").Value = 16105.92
$ws.Range("L132").Value = 15936.3999
$ws.Range("M132").Value = -13575.92
$ws.Range("N132").Value = -20996.3999

$ws.Range("H134").Value = 2785.3242
$ws.Range("I134").Value = 2773.5186
$ws.Range("J134").Value = 2817.2
$ws.Range("K134").Value = 8320.5558
$ws.Range("L134").Value = 8451.599999999999
$ws.Range("M134").Value = -5785.5558
$ws.Range("N134").Value = -13521.6

$ws.Range("H136").Value = 35715504
$ws.Range("I136").Value = 76923790
$ws.Range("J136").Value = 1652.2667
$ws.Range("K136").Value = 230771370
$ws.Range("L136").Value = 4956.800099999999
$ws.Range("M136").Value = -230768820
$ws.Range("N136").Value = -10056.8001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 4354.2
$ws.Range("I87").Value = 3455.25
$ws.Range("J87").Value = 7950
$ws.Range("K87").Value = 10365.75
$ws.Range("L87").Value = 23850
$ws.Range("M87").Value = -9117.75
$ws.Range("N87").Value = -26346

$ws.Range("H90").Value = 4354.2
$ws.Range("I90").Value = 3455.25
$ws.Range("J90").Value = 7950
$ws.Range("K90").Value = 31097.25
$ws.Range("L90").Value = 71550
$ws.Range("M90").Value = -24857.25
$ws.Range("N90").Value = -84030

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3975
$ws.Range("I70").Value = 3975
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 3975
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -3705
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 3975
$ws.Range("I73").Value = 3975
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 3975
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -3039
$ws.Range("N73").ClearContents()

$ws.Range("H80").Value = 4503
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 4503
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

$ws.Range("H102").Value = 4418.6
$ws.Range("I102").Value = 4146.875
$ws.Range("J102").Value = 5505.5
$ws.Range("K102").Value = 4146.875
$ws.Range("L102").Value = 5505.5
$ws.Range("M102").Value = -2524.875
$ws.Range("N102").Value = -8749.5

$ws.Range("H132").Value = 4283.863
$ws.Range("I132").Value = 4013
$ws.Range("J132").Value = 5739.75
$ws.Range("K132").Value = 12039
$ws.Range("L132").Value = 17219.25
$ws.Range("M132").Value = -9509
$ws.Range("N132").Value = -22279.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1908.7858
$ws.Range("I40").Value = 1451
$ws.Range("J40").Value = 2519.1667
$ws.Range("K40").Value = 1451
$ws.Range("L40").Value = 2519.1667
$ws.Range("M40").Value = -1315
$ws.Range("N40").Value = -2791.1667

$ws.Range("H68").Value = 1852.875
$ws.Range("I68").Value = 1761.5
$ws.Range("J68").Value = 1883.3334
$ws.Range("K68").Value = 1761.5
$ws.Range("L68").Value = 1883.3334
$ws.Range("M68").Value = -1012.5
$ws.Range("N68").Value = -3381.3334

$ws.Range("H71").Value = 1852.875
$ws.Range("I71").Value = 1761.5
$ws.Range("J71").Value = 1883.3334
$ws.Range("K71").Value = 8807.5
$ws.Range("L71").Value = 9416.666999999999
$ws.Range("M71").Value = -5063.5
$ws.Range("N71").Value = -16904.667

$ws.Range("H132").Value = 12900.641
$ws.Range("I132").Value = 5677.8066
$ws.Range("J132").Value = 40889.125
$ws.Range("K132").Value = 17033.4198
$ws.Range("L132").Value = 122667.375
$ws.Range("M132").Value = -14503.4198
$ws.Range("N132").Value = -127727.375

$ws.Range("H136").Value = 8894.1875
$ws.Range("I136").Value = 2476
$ws.Range("J136").Value = 15312.375
$ws.Range("K136").Value = 7428
$ws.Range("L136").Value = 45937.125
$ws.Range("M136").Value = -4878
$ws.Range("N136").Value = -51037.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3282.4285
$ws.Range("I132").Value = 3906.1064
$ws.Range("J132").Value = 2007.9565
$ws.Range("K132").Value = 11718.3192
$ws.Range("L132").Value = 6023.8695
$ws.Range("M132").Value = -9188.3192
$ws.Range("N132").Value = -11083.8695

$ws.Range("H136").Value = 7312.8823
$ws.Range("I136").Value = 8250.643
$ws.Range("K136").Value = 24751.929
$ws.Range("M136").Value = -24751.929
